# "CSV Variables.xlsx" - update the default Speed values (Speed Motorcycle /
# Speed Car / Speed Public transport) on the parameters sheet and leave the
# selection on the first of the edited cells, matching the addition of the
# new graphs/monitors that reference these speed inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.72
$ws.Range("C3").Value = 0.648
$ws.Range("C4").Value = 0.567

# Narrow the active selection from C2:C4 down to just C2.
$ws.Range("C2").Select()
